# Updated session 5 do files
# Row 2 and Row 3 are duplicate records; both receive the same updated
# values in columns K..BV (excluding a few columns that are unchanged:
# N, W, X, AB, AG, AH, AQ, AR, AV, BP).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "K"  = 0
    "L"  = 0
    "M"  = 1
    "O"  = 5.4117140769958496
    "P"  = 41.535858154296875
    "Q"  = 15.914735794067383
    "R"  = 25.621122360229492
    "S"  = 10.464022636413574
    "T"  = 37.738624572753906
    "U"  = 9.9870309829711914
    "V"  = 27.751594543457031
    "Y"  = 6.1416792869567871
    "Z"  = 30.753311157226563
    "AA" = 30.753311157226563
    "AC" = 17.944118499755859
    "AD" = 29.970684051513672
    "AE" = 5.6441683769226074
    "AF" = 24.326515197753906
    "AI" = 10.537778854370117
    "AJ" = 12.444417953491211
    "AK" = 7.4625377655029297
    "AL" = 4.9818801879882813
    "AM" = 19.420539855957031
    "AN" = 34.3587646484375
    "AO" = 23.82172966003418
    "AP" = 10.53703498840332
    "AS" = 17.804193496704102
    "AT" = 2.6125538349151611
    "AU" = 2.6125538349151611
    "AW" = 2.83695387840271
    "AX" = 31.2218017578125
    "AY" = 14.377996444702148
    "AZ" = 16.843805313110352
    "BA" = 1
    "BB" = 0
    "BC" = 4.3230438232421875
    "BD" = 9.7958612442016602
    "BE" = 9.7958612442016602
    "BF" = 0
    "BG" = 13.057753562927246
    "BH" = 5.6126728057861328
    "BI" = 5.6126728057861328
    "BJ" = 0
    "BK" = 1
    "BL" = 1
    "BM" = 18.734275817871094
    "BN" = 9.6759710311889648
    "BO" = 9.6759710311889648
    "BQ" = 14.227289199829102
    "BR" = 34.413803100585938
    "BS" = 20.556005477905273
    "BT" = 13.857797622680664
    "BU" = 2.328934907913208
    "BV" = 8.7675819396972656
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $ws.Range($col + "2").Value = $val
    $ws.Range($col + "3").Value = $val
}
